# Parts_List.xlsx restructuring:
#  - Split the single BOM table into two tables: "Parts We Don't Have" (existing
#    rows, minus the "Have"/"Yes" column and the Total row) and a new
#    "Parts We Have" table with one new part (RS-485 Transceivers).
#  - Add a hyperlink reference row for the new RS-485 transceiver part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- start clean: wipe existing values/formats/hyperlinks ---------------
$ws.Cells.Hyperlinks.Delete()
$ws.Cells.Clear()

$currencyFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# --- Table 1: "Parts We Don't Have" --------------------------------------
$ws.Range("A1").Value = "Parts We Don't Have"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Underline = $true
$ws.Range("A1").Font.Size = 14

$ws.Range("A2").Value = "Part"
$ws.Range("B2").Value = "Quantity"
$ws.Range("C2").Value = "Cost/Unit"
$ws.Range("D2").Value = "Cost"
$ws.Range("E2").Value = "Link"
$ws.Range("A2:E2").Font.Bold = $true
$ws.Range("A2:E2").Font.Size = 14

# Row 3 - 3.3V Voltage Reg
$ws.Range("A3").Value = "3.3V Voltage Reg"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.49
$ws.Range("E3").Value = "http://www.newark.com/stmicroelectronics/ld1117s33ctr/ic-ldo-volt-reg-3-3v-0-8a-sot/dp/89K0626?CMP=AFC-OP"

# Row 4 - W5100 TCP/IP Breakout
$ws.Range("A4").Value = "W5100 TCP/IP Breakout"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 24.95
$ws.Range("E4").Value = "https://www.sparkfun.com/products/9473"

# Row 5 - Dual 4:1 Muxes
$ws.Range("A5").Value = "Dual 4:1 Muxes"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0.57
$ws.Range("E5").Value = "http://www.newark.com/nxp/74hc4052d-653/ic-analog-mux-dmux-dual-4-x-1/dp/78R7402"

# Row 6 - GPS (no cost known yet)
$ws.Range("A6").Value = "GPS"
$ws.Range("B6").Value = 1
$ws.Range("E6").Value = "http://www.adafruit.com/products/746"

# Row 7 - Digital Temp (no quantity/cost known yet)
$ws.Range("A7").Value = "Digital Temp"
$ws.Range("E7").Value = "https://www.sparkfun.com/products/11295"

# Row 8 - Power Switch
$ws.Range("A8").Value = "Power Switch"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.75
$ws.Range("E8").Value = "https://www.sparkfun.com/products/9609"

# Row 9 - TI Stellaris EK-LM4F120XL
$ws.Range("A9").Value = "TI Stellaris EK-LM4F120XL"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 10.99
$ws.Range("E9").Value = "http://www.mouser.com/ProductDetail/Texas-Instruments/EK-LM4F120XL/?qs=t9Lg9qrXjEy2enepSwqR9A=="

# Cost formulas - D3 stands alone, D4:D9 share one formula group (D6/D7 blank)
$ws.Range("D3").Formula = "=(C3*B3)"
$ws.Range("D4:D9").Formula = "=(C4*B4)"
$ws.Range("D6").ClearContents()
$ws.Range("D7").ClearContents()

# Apply currency format to the whole cost columns (also re-adds the blank
# styled cells D6/D7/C6/C7 that ClearContents just emptied out)
$ws.Range("C3:C9").NumberFormat = $currencyFormat
$ws.Range("D3:D9").NumberFormat = $currencyFormat

# Left-align just the first part name in the table (mirrors the original
# file's centered A2, now left-aligned)
$ws.Range("A3").HorizontalAlignment = -4131

# Leftover formatted-but-empty remnants of the old "Total:" row
$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").Font.Size = 14
$ws.Range("C10").NumberFormat = $currencyFormat

# New reference link row
$ws.Range("E11").Value = "http://www.mouser.com/ProductDetail/Maxim-Integrated/MAX3077EESA+/?qs=sGAEpiMZZMuobhpKLk3hh6ov3TfCBqZhbNybjDy0atQ%3d"
$ws.Range("E11").Style = "Hyperlink"

# --- Table 2: "Parts We Have" --------------------------------------------
$ws.Range("A13").Value = "Parts We Have"
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").Font.Underline = $true
$ws.Range("A13").Font.Size = 14

$ws.Range("A14").Value = "Part"
$ws.Range("B14").Value = "Quantity"
$ws.Range("E14").Value = "Link"
$ws.Range("A14:E14").Font.Bold = $true
$ws.Range("A14:E14").Font.Size = 14

$ws.Range("A15").Value = "RS-485 Transceivers"
$ws.Range("B15").Value = 18
$ws.Range("E15").Value = "http://www.mouser.com/ProductDetail/Maxim-Integrated/MAX3077EESA+/?qs=sGAEpiMZZMuobhpKLk3hh6ov3TfCBqZhbNybjDy0atQ%3d"
$ws.Range("E15").Style = "Hyperlink"

$ws.Range("A17").Style = "Hyperlink"

# --- Hyperlinks -----------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E3"), "http://www.newark.com/stmicroelectronics/ld1117s33ctr/ic-ldo-volt-reg-3-3v-0-8a-sot/dp/89K0626?CMP=AFC-OP") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.sparkfun.com/products/9473") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "http://www.newark.com/nxp/74hc4052d-653/ic-analog-mux-dmux-dual-4-x-1/dp/78R7402") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"), "https://www.sparkfun.com/products/9609") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), "http://www.mouser.com/ProductDetail/Texas-Instruments/EK-LM4F120XL/?qs=t9Lg9qrXjEy2enepSwqR9A==") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E11"), "http://www.mouser.com/ProductDetail/Maxim-Integrated/MAX3077EESA+/?qs=sGAEpiMZZMuobhpKLk3hh6ov3TfCBqZhbNybjDy0atQ%3d") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E15"), "http://www.mouser.com/ProductDetail/Maxim-Integrated/MAX3077EESA+/?qs=sGAEpiMZZMuobhpKLk3hh6ov3TfCBqZhbNybjDy0atQ%3d") | Out-Null

# Note: GPS/Digital Temp hyperlinks (E6/E7) were not present in the original
# and are not in this diff either, keep them plain text.

# --- Column widths (character units) --------------------------------------
$ws.Columns("A:A").ColumnWidth = 24.17
$ws.Columns("B:B").ColumnWidth = 10.45
$ws.Columns("C:C").ColumnWidth = 11.45
$ws.Columns("D:D").ColumnWidth = 11.45
$ws.Columns("E:E").ColumnWidth = 102.59

# --- Selection / active cell ----------------------------------------------
$ws.Range("D17").Select()
